$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the loading_percent results for the "case with 380 kV" run (rows 2-25).
# Columns D, J, K, M, N keep their existing 0 values and are left untouched.
$cols = @("B","C","E","F","G","H","I","L","O")

$values = @(16.13996428727678, 11.68458322830679, 11.81306810524511, 16.86991607391245, 20.08359348480916, 12.25921542904854, 18.62252655626982, 9.644430804563072, 17.47800380851698)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $values[$i]
}

$values = @(15.36481217045498, 11.47277406483596, 11.87802688579383, 15.89584955866815, 20.28839276506119, 12.32947100404872, 18.78422403410015, 9.605792483442787, 17.61004981518622)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $values[$i]
}

$values = @(14.86746506304216, 11.34115157046004, 11.92048125042978, 15.26997757108489, 20.42508974919763, 12.37515964230774, 18.88854454662348, 9.583489560978663, 17.69641881102036)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $values[$i]
}

$values = @(14.65961587828911, 11.28717847264898, 11.93842805992277, 15.008197319934, 20.4835241531485, 12.39442016633994, 18.93232602511854, 9.574764799672472, 17.73294344951614)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $values[$i]
}

$values = @(14.62479672333347, 11.27819775765636, 11.94144715630453, 14.96433081551589, 20.49339128735222, 12.39765715282187, 18.93967270909193, 9.57333821813619, 17.73908851829904)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "6").Value = $values[$i]
}

$values = @(14.86468259826431, 11.3404249508783, 11.92072067018515, 15.26647399323133, 20.42586679678517, 12.37541679628318, 18.88912985232844, 9.583370414245492, 17.69690601856315)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "7").Value = $values[$i]
}

$values = @(15.87724633239141, 11.61191236064897, 11.83493273635117, 16.5399640634477, 20.15192116252763, 12.28291027241846, 18.67723626792698, 9.630817362275588, 17.5224335259444)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $values[$i]
}

$values = @(17.68566979439533, 12.12933556060565, 11.68708033432365, 19.00274580682531, 19.70268456699944, 12.12172889480456, 18.30153380138287, 9.734835036765702, 17.2223858832273)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value = $values[$i]
}

$values = @(18.89828374894878, 12.49717598148438, 11.59086245490442, 20.67494806633232, 19.42777593170594, 12.01560868081969, 18.0495673662913, 9.817519995809217, 17.02774448240386)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "10").Value = $values[$i]
}

$values = @(19.42350495588, 12.6612179016985, 11.54978275766214, 21.3917225636224, 19.31501967060603, 11.96999709050245, 17.94012101713592, 9.856397681700313, 16.94483365692677)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "11").Value = $values[$i]
}

$values = @(19.61851695818862, 12.72281461927013, 11.53461377200726, 21.65686569030329, 19.27411865008792, 11.95310783751096, 17.89941730308088, 9.871292661094946, 16.91425062838503)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "12").Value = $values[$i]
}

$values = @(19.57669123535358, 12.70957265060684, 11.53786346974628, 21.60004134736742, 19.28284702462385, 11.95672821297844, 17.90815065388278, 9.868077217940938, 16.92080099181578)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "13").Value = $values[$i]
}

$values = @(19.43962683991643, 12.66629624195243, 11.54852703900982, 21.4136618050453, 19.31161853413259, 11.968599928557, 17.9367574585536, 9.857619688466643, 16.94230124806439)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "14").Value = $values[$i]
}

$values = @(19.35516377871858, 12.63971874923673, 11.55510918724213, 21.29868154950795, 19.32947682794456, 11.97592155615749, 17.95437640074195, 9.851236390494353, 16.95557682589309)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "15").Value = $values[$i]
}

$values = @(18.86342072027399, 12.48638496843252, 11.59360125684754, 20.62722412089977, 19.43539477398827, 12.01864308518984, 18.05682382224062, 9.815003925199022, 17.03327654057601)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "16").Value = $values[$i]
}

$values = @(18.5549270894294, 12.39144093689248, 11.61790403580695, 20.20408069597325, 19.50354276521929, 12.04553335919718, 18.12099526489058, 9.793093942702534, 17.08238789991744)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "17").Value = $values[$i]
}

$values = @(18.37500831345933, 12.33652468271728, 11.63213556456316, 19.95656407809801, 19.54389580224786, 12.06125059213742, 18.15839219149893, 9.780611413165982, 17.11116539236846)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "18").Value = $values[$i]
}

$values = @(18.31366745426753, 12.31787980945837, 11.63699760145719, 19.87204792380568, 19.55775634812221, 12.06661523093053, 18.17113791424155, 9.776405836520203, 17.12099986853567)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "19").Value = $values[$i]
}

$values = @(18.58802421522954, 12.40158002593039, 11.61529075596827, 20.24955283636154, 19.49616843629746, 12.04264490206204, 18.11411370152654, 9.79541399226507, 17.07710504010159)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "20").Value = $values[$i]
}

$values = @(19.47999179119385, 12.67902212778764, 11.5453843857491, 21.46857628470577, 19.30311863752949, 11.96510252741592, 17.92833485183772, 9.860686699643725, 16.93596399432894)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "21").Value = $values[$i]
}

$values = @(20.04031356757576, 12.85728087592269, 11.50195237751784, 22.22866616901552, 19.18743898319693, 11.91665573654423, 17.81123688194148, 9.904349349947585, 16.8484636344694)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "22").Value = $values[$i]
}

$values = @(19.7433529829111, 12.76243684759354, 11.52492639202945, 21.82633154458858, 19.2482101430509, 11.94230851956687, 17.87333998615804, 9.880956987086961, 16.89472898010498)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "23").Value = $values[$i]
}

$values = @(18.57306897741545, 12.3969971730817, 11.61647141136003, 20.22900810905287, 19.49949871610141, 12.04394997093534, 18.11722328644357, 9.794364742679866, 17.07949173074828)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "24").Value = $values[$i]
}

$values = @(17.21633543071155, 11.99129947896736, 11.72489827842806, 18.34778573295695, 19.81464201656406, 12.16317080254456, 18.39893044632143, 9.705563329163519, 17.29903622442079)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "25").Value = $values[$i]
}
